# Fill in the "carrier" (D) column for the practice trials and the
# "pair_kind" (J) column for a couple of generic trials, plus populate the
# new unique_video / unique_audio trial blocks (rows 14-21) with their
# kind (C) and carrier (D) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Practice rows (2-5): carrier column D was empty, now filled in
$ws.Range("D2").Value = "can"
$ws.Range("D3").Value = "where"
$ws.Range("D4").Value = "do"
$ws.Range("D5").Value = "look"

# Generic rows (6-9): pair_kind column J was empty, now filled in
$ws.Range("J6").Value = "unique_video"
$ws.Range("J7").Value = "unique_video"
$ws.Range("J8").Value = "unique_audio"
$ws.Range("J9").Value = "unique_audio"

# New unique_video / unique_audio trial rows (14-21): kind (C) and carrier (D)
$ws.Range("C14").Value = "unique_video"
$ws.Range("D14").Value = "look"

$ws.Range("C15").Value = "unique_video"
$ws.Range("D15").Value = "look"

$ws.Range("C16").Value = "unique_video"
$ws.Range("D16").Value = "where"

$ws.Range("C17").Value = "unique_video"
$ws.Range("D17").Value = "where"

$ws.Range("C18").Value = "unique_audio"
$ws.Range("D18").Value = "can"

$ws.Range("C19").Value = "unique_audio"
$ws.Range("D19").Value = "can"

$ws.Range("C20").Value = "unique_audio"
$ws.Range("D20").Value = "do"

$ws.Range("C21").Value = "unique_audio"
$ws.Range("D21").Value = "do"
